$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-6 got reassigned to new positions (e.g. re-sorted by date)
# while every column within a given record stayed together. Capture each
# full row (columns A:T) before writing anything back, then re-assign rows
# to their new positions according to the observed permutation:
#   new row 2 <- old row 5
#   new row 3 <- old row 6
#   new row 4 <- old row 3
#   new row 5 <- old row 2
#   new row 6 <- old row 4

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$rows = @{}
for ($r = 2; $r -le 6; $r++) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$r").Value2
    }
    $rows[$r] = $rowData
}

$mapping = @{
    2 = 5
    3 = 6
    4 = 3
    5 = 2
    6 = 4
}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $rowData = $rows[$oldRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value2 = $rowData[$col]
    }
}
